$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("study")
$ws.Rows.Item(7).Insert()
$ws.Cells.Item(7, 1).Value = "businessTherapeuticAreas"
